$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 14182.777
$ws.Range("I86").Value = 6151.5
$ws.Range("J86").Value = 16477.428
$ws.Range("K86").Value = 6151.5
$ws.Range("L86").Value = 16477.428
$ws.Range("M86").Value = -5028.5
$ws.Range("N86").Value = -18723.428

$ws.Range("H89").Value = 14182.777
$ws.Range("I89").Value = 6151.5
$ws.Range("J89").Value = 16477.428
$ws.Range("K89").Value = 30757.5
$ws.Range("L89").Value = 82387.14
$ws.Range("M89").Value = -25141.5
$ws.Range("N89").Value = -93619.14

$ws.Range("H92").Value = 100000750
$ws.Range("I92").Value = 100000750
$ws.Range("K92").Value = 100000750
$ws.Range("M92").Value = -99999502

$ws.Range("H112").Value = 1134.4
$ws.Range("J112").Value = 1134.4
$ws.Range("L112").Value = 3403.2
$ws.Range("N112").Value = -5619.200000000001

$ws.Range("H116").Value = 5928
$ws.Range("J116").Value = 5928
$ws.Range("L116").Value = 5928
$ws.Range("N116").Value = -12812

$ws.Range("H129").Value = 193874.22
$ws.Range("J129").Value = 197665.86
$ws.Range("L129").Value = 592997.58
$ws.Range("N129").Value = -602997.58

$ws.Range("H132").Value = 2356.6135
$ws.Range("I132").Value = 2494.9
$ws.Range("K132").Value = 7484.700000000001
$ws.Range("M132").Value = -4954.700000000001

$ws.Range("H138").Value = 1849.0834
$ws.Range("I138").Value = 797.3182
$ws.Range("J138").Value = 2222.2903
$ws.Range("K138").Value = 2391.9546
$ws.Range("L138").Value = 6666.8709
$ws.Range("M138").Value = 2748.0454
$ws.Range("N138").Value = -16946.8709

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1504.721
$ws.Range("I2").Value = 1447.2059
$ws.Range("K2").Value = 1447.2059
$ws.Range("M2").Value = -1334.2059

$ws.Range("H32").Value = 19325.2
$ws.Range("I32").Value = 21892.41
$ws.Range("K32").Value = 21892.41
$ws.Range("M32").Value = -21605.41

$ws.Range("H45").Value = 2955.611
$ws.Range("I45").Value = 2784.7856
$ws.Range("K45").Value = 2784.7856
$ws.Range("M45").Value = -2407.7856

$ws.Range("H74").Value = 76924776
$ws.Range("I74").Value = 125000696
$ws.Range("J74").Value = 3311
$ws.Range("K74").Value = 125000696
$ws.Range("L74").Value = 3311
$ws.Range("M74").Value = -124999822
$ws.Range("N74").Value = -5059

$ws.Range("H77").Value = 76924776
$ws.Range("I77").Value = 125000696
$ws.Range("J77").Value = 3311
$ws.Range("K77").Value = 625003480
$ws.Range("L77").Value = 16555
$ws.Range("M77").Value = -624999112
$ws.Range("N77").Value = -25291

$ws.Range("H97").Value = 1240.52
$ws.Range("I97").Value = 1435.9412
$ws.Range("K97").Value = 1435.9412
$ws.Range("M97").Value = -939.9412

$ws.Range("H110").Value = 599.8333
$ws.Range("I110").Value = 400
$ws.Range("J110").Value = 639.8
$ws.Range("K110").Value = 400
$ws.Range("L110").Value = 639.8
$ws.Range("M110").Value = 1645
$ws.Range("N110").Value = -4729.8

$ws.Range("H116").Value = 1504.721
$ws.Range("I116").Value = 1447.2059
$ws.Range("K116").Value = 1447.2059
$ws.Range("M116").Value = 846.7941000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1504.721
$ws.Range("I3").Value = 1447.2059
$ws.Range("K3").Value = 1447.2059
$ws.Range("M3").Value = -1333.2059

$ws.Range("H94").Value = 1748.9333
$ws.Range("I94").Value = 851.1667
$ws.Range("J94").Value = 5340
$ws.Range("K94").Value = 851.1667
$ws.Range("L94").Value = 5340
$ws.Range("M94").Value = -400.1667
$ws.Range("N94").Value = -6242

$ws.Range("H107").Value = 1392.5333
$ws.Range("I107").Value = 1023
$ws.Range("K107").Value = 1023
$ws.Range("M107").Value = 897

$ws.Range("H134").Value = 48453.957
$ws.Range("I134").Value = 69021.44
$ws.Range("J134").Value = 1442.5714
$ws.Range("K134").Value = 207064.32
$ws.Range("L134").Value = 4327.7142
$ws.Range("M134").Value = -204529.32
$ws.Range("N134").Value = -9397.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 10381.77
$ws.Range("I31").Value = 18014.947
$ws.Range("J31").Value = 3130.25
$ws.Range("K31").Value = 18014.947
$ws.Range("L31").Value = 3130.25
$ws.Range("M31").Value = -17719.947
$ws.Range("N31").Value = -3720.25

$ws.Range("H34").Value = 10381.77
$ws.Range("I34").Value = 18014.947
$ws.Range("J34").Value = 3130.25
$ws.Range("K34").Value = 18014.947
$ws.Range("L34").Value = 3130.25
$ws.Range("M34").Value = -17812.947
$ws.Range("N34").Value = -3534.25

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H60").Value = 14495
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 14495
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 14495
$ws.Range("N60").Value = -15517
$ws.Range("M60").ClearContents()

$ws.Range("H94").Value = 2552.8948
$ws.Range("I94").Value = 909.3333
$ws.Range("J94").Value = 3311.4614
$ws.Range("K94").Value = 909.3333
$ws.Range("L94").Value = 3311.4614
$ws.Range("M94").Value = -458.3333
$ws.Range("N94").Value = -4213.4614

$ws.Range("H107").Value = 1205.8096
$ws.Range("I107").Value = 318.25
$ws.Range("J107").Value = 1752
$ws.Range("K107").Value = 318.25
$ws.Range("L107").Value = 1752
$ws.Range("M107").Value = 1601.75
$ws.Range("N107").Value = -5592

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 20825.5
$ws.Range("I132").Value = 28614.79
$ws.Range("J132").Value = 4381.4443
$ws.Range("K132").Value = 85844.37
$ws.Range("L132").Value = 13144.3329
$ws.Range("M132").Value = -83314.37
$ws.Range("N132").Value = -18204.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1050
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1225
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3675
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -5297

$ws.Range("H71").Value = 1050
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1225
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 11025
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -19137

$ws.Range("H120").Value = 16012
$ws.Range("I120").Value = 10030
$ws.Range("K120").Value = 30090
$ws.Range("M120").Value = -25252

$ws.Range("H131").Value = 766.67
$ws.Range("J131").Value = 806.89246
$ws.Range("L131").Value = 2420.67738
$ws.Range("N131").Value = -12500.67738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 718.46155
$ws.Range("I97").Value = 718.46155
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 718.46155
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -222.46155
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 27779048
$ws.Range("I102").Value = 31251210
$ws.Range("K102").Value = 31251210
$ws.Range("M102").Value = -31249588

$ws.Range("H104").Value = 29999
$ws.Range("J104").Value = 29999
$ws.Range("L104").Value = 29999
$ws.Range("N104").Value = -36987

$ws.Range("H122").Value = 49384456
$ws.Range("I122").Value = 18519876
$ws.Range("J122").Value = 111113610
$ws.Range("K122").Value = 55559628
$ws.Range("L122").Value = 333340830
$ws.Range("M122").Value = -55557178
$ws.Range("N122").Value = -333345730

$ws.Range("H126").Value = 5206.875
$ws.Range("I126").Value = 3956.1875
$ws.Range("K126").Value = 11868.5625
$ws.Range("M126").Value = -9398.5625

$ws.Range("H132").Value = 75627.664
$ws.Range("I132").Value = 89015.414
$ws.Range("K132").Value = 267046.242
$ws.Range("M132").Value = -264516.242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 27698.75
$ws.Range("J59").Value = 27698.75
$ws.Range("L59").Value = 27698.75
$ws.Range("N59").Value = -29006.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -15040

$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -627
$ws.Range("N96").ClearContents()

$ws.Range("H107").Value = 3247469.8
$ws.Range("I107").Value = 805.7778
$ws.Range("J107").Value = 9091465
$ws.Range("K107").Value = 2417.3334
$ws.Range("L107").Value = 27274395
$ws.Range("M107").Value = -497.3334
$ws.Range("N107").Value = -27278235

$ws.Range("H126").Value = 2699.75
$ws.Range("I126").Value = 1459.8
$ws.Range("J126").Value = 4766.3335
$ws.Range("K126").Value = 4379.4
$ws.Range("L126").Value = 14299.0005
$ws.Range("M126").Value = -1909.4
$ws.Range("N126").Value = -19239.0005

$ws.Range("H132").Value = 1695.4546
$ws.Range("I132").Value = 919
$ws.Range("J132").Value = 3766
$ws.Range("K132").Value = 2757
$ws.Range("L132").Value = 11298
$ws.Range("M132").Value = -227
$ws.Range("N132").Value = -16358
